# CIERRE 14 MAY 22
# Roll the payroll receipts sheet (Hoja1) forward from SEMANA 18
# (02-Al-01 MAYO 2022) to SEMANA 19 (09-Al-15 MAYO 2022), updating the
# week-header text and the two pay figures that changed for the new
# closing, then leave the selection where the user ended up (H41).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Week header (B9) drives H9 (=B9), B27 (=B9), H27 (=B27) and B43 (=H27)
# via formulas, so a single edit here ripples through all of them.
$ws.Range("B9").Value = "SEMANA   19  DEL    09      Al   15   DE   MAYO          2022"

# Updated pay figures for the new closing week.
$ws.Range("K21").Value = 1120
$ws.Range("E40").Value = 1250

# Leave the cursor/selection on H41, matching the saved view state.
$ws.Range("H41").Select() | Out-Null
